$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs", "Col4a3", "Cd47", "ECs", 3, 1, 0.6813643333333333, 2.044093, 0.3450055461675409, 0.3450055461675409, 3, 1, 46.33695966666667, 139.010879, 0.1993490803952133, 0.1993490803952133, 31.57235163197189, 284.151164687747, 0.06877653835974759, 0.06877653835974759),
    @("ECs", "Col4a3", "Cd47", "FAPs", 3, 1, 0.6813643333333333, 2.044093, 0.3450055461675409, 0.3450055461675409, 3, 1, 84.50960033333332, 253.528801, 0.3635739425333109, 0.3635739425333109, 57.5818274913881, 518.236447422493, 0.125435026615991, 0.1254350266159911),
    @("ECs", "Col4a3", "Cd47", "M2", 3, 1, 0.6813643333333333, 2.044093, 0.3450055461675409, 0.3450055461675409, 3, 1, 72.52790466666666, 217.583714, 0.3120267536390091, 0.3120267536390091, 49.41792741126688, 444.761346701402, 0.107650960558111, 0.1076509605581111),
    @("ECs", "Col4a3", "Cd47", "sCs", 3, 1, 0.6813643333333333, 2.044093, 0.3450055461675409, 0.3450055461675409, 3, 1, 29.06683666666666, 87.20050999999998, 0.1250502234324667, 0.1250502234324667, 19.80510578749222, 178.24595208743, 0.04314302063369121, 0.04314302063369121),
    @("FAPs", "Col4a3", "Cd47", "ECs", 3, 1, 0.845217, 2.535651, 0.4279715542029013, 0.4279715542029013, 3, 1, 46.33695966666667, 139.010879, 0.1993490803952133, 0.1993490803952133, 39.164786038581, 352.483074347229, 0.08531573576565857, 0.08531573576565857),
    @("FAPs", "Col4a3", "Cd47", "FAPs", 3, 1, 0.845217, 2.535651, 0.4279715542029013, 0.4279715542029013, 3, 1, 84.50960033333332, 253.528801, 0.3635739425333109, 0.3635739425333109, 71.428950864939, 642.860557784451, 0.1555993052536574, 0.1555993052536574),
    @("FAPs", "Col4a3", "Cd47", "M2", 3, 1, 0.845217, 2.535651, 0.4279715542029013, 0.4279715542029013, 3, 1, 72.52790466666666, 217.583714, 0.3120267536390091, 0.3120267536390091, 61.30181799864599, 551.716361987814, 0.1335385747077725, 0.1335385747077725),
    @("FAPs", "Col4a3", "Cd47", "sCs", 3, 1, 0.845217, 2.535651, 0.4279715542029013, 0.4279715542029013, 3, 1, 29.06683666666666, 87.20050999999998, 0.1250502234324667, 0.1250502234324667, 24.56778448688999, 221.1100603820099, 0.05351793847581286, 0.05351793847581286),
    @("M2", "Col4a3", "Cd47", "ECs", 1, 0.3333333333333333, 0.009069333333333334, 0.027208, 0.004592213221280271, 0.004592213221280271, 3, 1, 46.33695966666667, 139.010879, 0.1993490803952133, 0.1993490803952133, 0.4202453328702223, 3.782207995832, 0.0009154534826409621, 0.000915453482640962),
    @("M2", "Col4a3", "Cd47", "FAPs", 1, 0.3333333333333333, 0.009069333333333334, 0.027208, 0.004592213221280271, 0.004592213221280271, 3, 1, 84.50960033333332, 253.528801, 0.3635739425333109, 0.3635739425333109, 0.7664457352897778, 6.898011617608, 0.001669609065814464, 0.001669609065814464),
    @("M2", "Col4a3", "Cd47", "M2", 1, 0.3333333333333333, 0.009069333333333334, 0.027208, 0.004592213221280271, 0.004592213221280271, 3, 1, 72.52790466666666, 217.583714, 0.3120267536390091, 0.3120267536390091, 0.6577797433902222, 5.920017690511999, 0.001432893383454219, 0.001432893383454219),
    @("M2", "Col4a3", "Cd47", "sCs", 1, 0.3333333333333333, 0.009069333333333334, 0.027208, 0.004592213221280271, 0.004592213221280271, 3, 1, 29.06683666666666, 87.20050999999998, 0.1250502234324667, 0.1250502234324667, 0.2636168306755555, 2.372551476079999, 0.0005742572893706257, 0.0005742572893706257),
    @("sCs", "Col4a3", "Cd47", "ECs", 3, 1, 0.4392866666666667, 1.31786, 0.2224306864082776, 0.2224306864082776, 3, 1, 46.33695966666667, 139.010879, 0.1993490803952133, 0.1993490803952133, 20.35520855543778, 183.19687699894, 0.04434135278716621, 0.04434135278716621),
    @("sCs", "Col4a3", "Cd47", "FAPs", 3, 1, 0.4392866666666667, 1.31786, 0.2224306864082776, 0.2224306864082776, 3, 1, 84.50960033333332, 253.528801, 0.3635739425333109, 0.3635739425333109, 37.12394063176222, 334.11546568586, 0.08087000159784802, 0.08087000159784802),
    @("sCs", "Col4a3", "Cd47", "M2", 3, 1, 0.4392866666666667, 1.31786, 0.2224306864082776, 0.2224306864082776, 3, 1, 72.52790466666666, 217.583714, 0.3120267536390091, 0.3120267536390091, 31.86054148133777, 286.74487333204, 0.06940432498967132, 0.06940432498967132),
    @("sCs", "Col4a3", "Cd47", "sCs", 3, 1, 0.4392866666666667, 1.31786, 0.2224306864082776, 0.2224306864082776, 3, 1, 29.06683666666666, 87.20050999999998, 0.1250502234324667, 0.1250502234324667, 12.76867378984444, 114.9180641086, 0.02781500703359206, 0.02781500703359206)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
